$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update all data rows (2-10) with the refreshed scrape content ---
# Rows 2-8 reuse their existing cell positions; the F-column hyperlink
# relationships already anchored at F2:F8 are left untouched (same as source).
# Rows 9-10 are newly appended and get brand-new hyperlinks below.

# Row 2
$ws.Cells.Item(2,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(2,2).Value = "【完全在宅】AI×Web開発エンジニア募集!業務自動化・AI機能開発"
$ws.Cells.Item(2,3).Value = "システム開発"
$ws.Cells.Item(2,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(2,5).Value = "期限情報なし"
$ws.Cells.Item(2,6).Value = "https://www.lancers.jp/work/detail/5483480"
$ws.Cells.Item(2,7).Value = 435
$ws.Cells.Item(2,8).Value = "🔥AI,Ai ◆開発,自動化"

# Row 3
$ws.Cells.Item(3,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(3,2).Value = "【時給3000~4000円以上/フルリモート】AI駆動開発でのSaaS開発の開発パートナー"
$ws.Cells.Item(3,3).Value = "システム開発"
$ws.Cells.Item(3,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3,5).Value = "期限情報なし"
$ws.Cells.Item(3,6).Value = "https://www.lancers.jp/work/detail/5483313"
$ws.Cells.Item(3,7).Value = 368
$ws.Cells.Item(3,8).Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Cells.Item(4,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(4,2).Value = "【急募】ウェブサイトのAIチャットサポートの実装とLINE公式アカウントのAIチャットサポート"
$ws.Cells.Item(4,3).Value = "システム開発"
$ws.Cells.Item(4,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(4,5).Value = "期限情報なし"
$ws.Cells.Item(4,6).Value = "https://www.lancers.jp/work/detail/5483343"
$ws.Cells.Item(4,7).Value = 338
$ws.Cells.Item(4,8).Value = "🔥AI,Ai ◇サイト"

# Row 5
$ws.Cells.Item(5,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(5,2).Value = "【急募】人事評価・賞与計算を自動化する社内向けWebシステム開発"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5483345"
$ws.Cells.Item(5,7).Value = 203
$ws.Cells.Item(5,8).Value = "◆開発,システム開発"

# Row 6
$ws.Cells.Item(6,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(6,2).Value = "本人顔ベースのリアルタイム顔変換システム開発"
$ws.Cells.Item(6,3).Value = "システム開発"
$ws.Cells.Item(6,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(6,5).Value = "期限情報なし"
$ws.Cells.Item(6,6).Value = "https://www.lancers.jp/work/detail/5483207"
$ws.Cells.Item(6,7).Value = 125
$ws.Cells.Item(6,8).Value = "◆開発,システム開発"

# Row 7
$ws.Cells.Item(7,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(7,2).Value = "【急募】新しいWebサービスの開発パートナーを探しています!"
$ws.Cells.Item(7,3).Value = "システム開発"
$ws.Cells.Item(7,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(7,5).Value = "期限情報なし"
$ws.Cells.Item(7,6).Value = "https://www.lancers.jp/work/detail/5483482"
$ws.Cells.Item(7,7).Value = 75
$ws.Cells.Item(7,8).Value = "◆開発"

# Row 8
$ws.Cells.Item(8,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(8,2).Value = "【PM/フルスタックエンジニア】新規SaaS開発のパートナー募集"
$ws.Cells.Item(8,3).Value = "システム開発"
$ws.Cells.Item(8,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(8,5).Value = "期限情報なし"
$ws.Cells.Item(8,6).Value = "https://www.lancers.jp/work/detail/5483306"
$ws.Cells.Item(8,7).Value = 75
$ws.Cells.Item(8,8).Value = "◆開発"

# Row 9
$ws.Cells.Item(9,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(9,2).Value = "カフェ掲載プラットフォーム「チャヤドコ」開発(要件定義~ベータ版リリース)"
$ws.Cells.Item(9,3).Value = "システム開発"
$ws.Cells.Item(9,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(9,5).Value = "期限情報なし"
$ws.Cells.Item(9,6).Value = "https://www.lancers.jp/work/detail/5483311"
$ws.Cells.Item(9,7).Value = 68
$ws.Cells.Item(9,8).Value = "◆開発"

# Row 10
$ws.Cells.Item(10,1).Value = "2026-02-01 18:31:42"
$ws.Cells.Item(10,2).Value = "【月次継続】hacomono・Stripeデータ管理のプロを求む!"
$ws.Cells.Item(10,3).Value = "システム開発"
$ws.Cells.Item(10,4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10,5).Value = "期限情報なし"
$ws.Cells.Item(10,6).Value = "https://www.lancers.jp/work/detail/5483295"
$ws.Cells.Item(10,7).Value = 38
$ws.Cells.Item(10,8).Value = "◇管理"

# --- New hyperlinks for the two newly appended rows ---
$ws.Hyperlinks.Add($ws.Cells.Item(9,6), "https://www.lancers.jp/work/detail/5483311") | Out-Null
$ws.Cells.Item(9,6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(10,6), "https://www.lancers.jp/work/detail/5483295") | Out-Null
$ws.Cells.Item(10,6).Style = "Hyperlink"

# --- Column H got wider to fit the longer skill-summary text ---
$ws.Columns.Item(8).ColumnWidth = 15.1666666667